$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 196; this shifts existing rows 196-296 down
# to 197-297 and extends the used range from A1:R296 to A1:R297.
$ws.Rows(196).Insert()

# Populate the newly inserted row 196 with the new record.
$ws.Range("A196").Value = 11
$ws.Range("B196").Value = "Vega Monumental Concepción"
$ws.Range("C196").Value = "Bíobío"
$ws.Range("D196").Value = 44777
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = 100114001
$ws.Range("G196").Value = "Papa"
$ws.Range("H196").Value = "Patagonia"
$ws.Range("I196").Value = "1a (guarda)"
$ws.Range("J196").Value = 220
$ws.Range("K196").Value = 7000
$ws.Range("L196").Value = 7500
$ws.Range("M196").Value = 7227
$ws.Range("N196").Value = "`$/saco 25 kilos"
$ws.Range("O196").Value = "Región de Los Lagos"
$ws.Range("P196").Value = 289
$ws.Range("Q196").Value = 25
$ws.Range("R196").Value = "Hortaliza"
